# Update Sheets via scheduled runner: refresh market price data across ALC, ARM, BSM, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2899.5
$ws.Range("I29").Value = 2899.5
$ws.Range("K29").Value = 8698.5
$ws.Range("M29").Value = -8417.5
$ws.Range("H76").Value = 3238.4614
$ws.Range("I76").Value = 3154.5454
$ws.Range("J76").Value = 3700
$ws.Range("K76").Value = 3154.5454
$ws.Range("L76").Value = 3700
$ws.Range("M76").Value = -2839.5454
$ws.Range("N76").Value = -4330
$ws.Range("H79").Value = 3238.4614
$ws.Range("I79").Value = 3154.5454
$ws.Range("J79").Value = 3700
$ws.Range("K79").Value = 3154.5454
$ws.Range("L79").Value = 3700
$ws.Range("M79").Value = -2062.5454
$ws.Range("N79").Value = -5884
$ws.Range("H141").Value = 581789.4399999999
$ws.Range("I141").Value = 1748.5333
$ws.Range("J141").Value = 1669366.1
$ws.Range("K141").Value = 5245.5999
$ws.Range("L141").Value = 5008098.300000001
$ws.Range("M141").Value = -65.59990000000016
$ws.Range("N141").Value = -5018458.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1781.0646
$ws.Range("I45").Value = 1083.96
$ws.Range("J45").Value = 4685.6665
$ws.Range("K45").Value = 1083.96
$ws.Range("L45").Value = 4685.6665
$ws.Range("M45").Value = -706.96
$ws.Range("N45").Value = -5439.6665
$ws.Range("H97").Value = 344.3684
$ws.Range("I97").Value = 344.3684
$ws.Range("K97").Value = 344.3684
$ws.Range("M97").Value = 151.6316
$ws.Range("H122").Value = 4119.933
$ws.Range("I122").Value = 2571.4285
$ws.Range("J122").Value = 5474.875
$ws.Range("K122").Value = 7714.2855
$ws.Range("L122").Value = 16424.625
$ws.Range("M122").Value = -5264.2855
$ws.Range("N122").Value = -21324.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1623.9565
$ws.Range("I99").Value = 1045.7142
$ws.Range("J99").Value = 2523.4443
$ws.Range("K99").Value = 1045.7142
$ws.Range("L99").Value = 2523.4443
$ws.Range("M99").Value = 452.2858000000001
$ws.Range("N99").Value = -5519.4443
$ws.Range("H134").Value = 3483.2778
$ws.Range("I134").Value = 2446.6
$ws.Range("J134").Value = 8666.666999999999
$ws.Range("K134").Value = 7339.799999999999
$ws.Range("L134").Value = 26000.001
$ws.Range("M134").Value = -4804.799999999999
$ws.Range("N134").Value = -31070.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 74.64286
$ws.Range("I2").Value = 84.25
$ws.Range("J2").Value = 61.833332
$ws.Range("K2").Value = 505.5
$ws.Range("L2").Value = 370.999992
$ws.Range("M2").Value = -392.5
$ws.Range("N2").Value = -596.999992
$ws.Range("H17").Value = 969.75
$ws.Range("I17").Value = 440
$ws.Range("J17").Value = 1499.5
$ws.Range("K17").Value = 1320
$ws.Range("L17").Value = 4498.5
$ws.Range("M17").Value = -1151
$ws.Range("N17").Value = -4836.5
$ws.Range("H34").Value = 18350
$ws.Range("J34").Value = 21980
$ws.Range("L34").Value = 65940
$ws.Range("N34").Value = -66108
$ws.Range("H39").Value = 3685.5715
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3685.5715
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = 11056.7145
$ws.Range("N39").Value = -11644.7145
$ws.Range("L39").ClearContents()
$ws.Range("H50").Value = 86163.414
$ws.Range("I50").Value = 52.5
$ws.Range("J50").Value = 103385.6
$ws.Range("K50").Value = 157.5
$ws.Range("L50").Value = 310156.8
$ws.Range("M50").Value = 323.5
$ws.Range("N50").Value = -311118.8
$ws.Range("H53").Value = 86163.414
$ws.Range("I53").Value = 52.5
$ws.Range("J53").Value = 103385.6
$ws.Range("K53").Value = 157.5
$ws.Range("L53").Value = 310156.8
$ws.Range("M53").Value = 323.5
$ws.Range("N53").Value = -311118.8
$ws.Range("H55").Value = 2559.2307
$ws.Range("I55").Value = 490
$ws.Range("J55").Value = 2935.4546
$ws.Range("K55").Value = 1470
$ws.Range("L55").Value = 8806.363799999999
$ws.Range("M55").Value = -1293
$ws.Range("N55").Value = -9160.363799999999
$ws.Range("H105").Value = 3642.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3743.6365
$ws.Range("I70").Value = 3688.889
$ws.Range("K70").Value = 3688.889
$ws.Range("M70").Value = -3418.889
$ws.Range("H73").Value = 3743.6365
$ws.Range("I73").Value = 3688.889
$ws.Range("K73").Value = 3688.889
$ws.Range("M73").Value = -2752.889
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H122").Value = 3150
$ws.Range("J122").Value = 5333.3335
$ws.Range("L122").Value = 16000.0005
$ws.Range("N122").Value = -20900.0005
$ws.Range("H132").Value = 3125.8965
$ws.Range("I132").Value = 2161.2856
$ws.Range("J132").Value = 4026.2
$ws.Range("K132").Value = 6483.8568
$ws.Range("L132").Value = 12078.6
$ws.Range("M132").Value = -3953.8568
$ws.Range("N132").Value = -17138.6
$ws.Range("H140").Value = 38333.332
$ws.Range("J140").Value = 38333.332
$ws.Range("L140").Value = 38333.332
$ws.Range("N140").Value = -48693.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 35102.5
$ws.Range("I7").Value = 200
$ws.Range("J7").Value = 70005
$ws.Range("K7").Value = 200
$ws.Range("L7").Value = 70005
$ws.Range("M7").Value = -87
$ws.Range("N7").Value = -70231
$ws.Range("H46").Value = 36222.57
$ws.Range("J46").Value = 36222.57
$ws.Range("L46").Value = 36222.57
$ws.Range("N46").Value = -36684.57
$ws.Range("H81").Value = 847.6667
$ws.Range("I81").Value = 662.25
$ws.Range("J81").Value = 996
$ws.Range("K81").Value = 1324.5
$ws.Range("L81").Value = 1992
$ws.Range("M81").Value = -263.5
$ws.Range("N81").Value = -4114
$ws.Range("H84").Value = 847.6667
$ws.Range("I84").Value = 662.25
$ws.Range("J84").Value = 996
$ws.Range("K84").Value = 6622.5
$ws.Range("L84").Value = 9960
$ws.Range("M84").Value = -1318.5
$ws.Range("N84").Value = -20568
$ws.Range("H132").Value = 13284.5
$ws.Range("I132").Value = 1986
$ws.Range("K132").Value = 5958
$ws.Range("M132").Value = -3428
$ws.Range("H134").Value = 36222.57
$ws.Range("J134").Value = 36222.57
$ws.Range("L134").Value = 108667.71
$ws.Range("N134").Value = -113737.71
